# Updated database excel download
#
# "Trees" sheet restructuring:
#   - old column A "Name"       -> becomes new column B "Tree ID" (keeps the
#     original tree-name values)
#   - a new column A "User" is inserted, with every data row set to the
#     constant value "tang"
#   - old column C "Stem Count" is removed entirely
#   - all other columns (Circumf, Tapping Date, Tap Height, Latitude,
#     Longitude, Start/End of Season Notes) keep their data, just shifted
#     right by one column because of the new "User" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trees")

# Locate the last used row/column on the Trees sheet.
$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count

# 1) Remember the existing tree names (old column A, rows below the header).
$treeIds = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $treeIds += $ws.Cells.Item($r, 1).Text
}

# 2) Insert a new blank column before column B - this pushes the old
#    Circumf/Stem Count/... columns one to the right and leaves room for the
#    "Tree ID" column right after "User" (old column A).
$ws.Columns.Item(2).Insert()

# 3) Re-label the headers: old A1 "Name" -> "User"; new B1 -> "Tree ID".
$ws.Cells.Item(1, 1).Value = "User"
$ws.Cells.Item(1, 2).Value = "Tree ID"

# 4) Populate the new "Tree ID" column with the old tree names, and set the
#    "User" column to the constant "tang" for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $treeIds[$r - 2]
    $ws.Cells.Item($r, 1).Value = "tang"
}

# 5) Delete the old "Stem Count" column, which is now column D (it shifted
#    right by one when the "Tree ID" column was inserted in step 2).
$ws.Columns.Item(4).Delete()
